$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: remove the trailing empty placeholder cells AX24:AZ24 entirely
# (these were blank inline-string placeholder cells; clearing removes the cell
# node from the sheet XML, matching the target which no longer has them).
$ws.Range("AX24:AZ24").ClearContents()

# Row 25: new student record appended by the "add aluno" form submission.
# Force text storage (matches source data where every field -- including
# numeric-looking codes and ISO dates -- is stored as literal text) by setting
# NumberFormat to Text ("@") before writing each value, then drop the format
# back off so no stray style lingers on the cell.
$row = 25
$values = [ordered]@{
    "A" = "1"
    "B" = "q"
    "C" = "1"
    "D" = "Branca"
    "E" = "Masculino"
    "F" = "1"
    "G" = "1"
    "H" = "a"
    "I" = "a"
    "J" = "a"
    "K" = "2000-01-01"
    "L" = "a"
    "M" = "aa"
    "N" = "1"
    "O" = "NÃO"
    "P" = "NÃO"
    "Q" = "NÃO"
    "R" = "NÃO"
    "S" = "NÃO"
    "T" = "NÃO"
    "U" = "NÃO"
    "V" = "NÃO"
    "W" = "NÃO"
    "X" = "NÃO"
    "Y" = "NÃO"
    "Z" = "NÃO"
    "AA" = "NÃO"
    "AB" = "NÃO"
    "AC" = "a"
    "AD" = "1"
    "AE" = "1"
    "AF" = "1"
    "AG" = "1"
    "AH" = "Rural"
    "AI" = "1"
    "AJ" = "1"
    "AK" = "a"
    "AL" = "a"
    "AM" = "a"
    "AN" = "1"
    "AO" = "1/1/2000"
    "AP" = "1"
    "AQ" = "1/1/2000"
    "AR" = "Manhã"
    "AS" = "01. Berçário I"
    "AT" = "02 - Escola Municipal"
    "AU" = "NÃO"
    "AV" = "NÃO"
    "AW" = "1"
    "AX" = ""
    "AY" = ""
    "AZ" = ""
    "BA" = "2000-01-01"
}

foreach ($col in $values.Keys) {
    $cell = $ws.Range("$col$row")
    $cell.NumberFormat = "@"
    $cell.Value = $values[$col]
    $cell.ClearFormats()
}

